# Generate Report for Handback
# -----------------------------------------------------------------------
# The German (de-de) and Chinese (zh-cn) handback round finished: every
# row's status flips from "Ready for handoff" to "Handed back: in sync
# with en-US", and the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns get populated for both locales.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column width, expressed via the COM ColumnWidth property, is quantized
# to Excel's internal pixel grid -- these inputs land in the same
# rounding bucket as the authored target widths (~29.98 -> 30, ~40 -> 40).
$wideColWidth = 29.15
$maxColWidth = 39.15

# ----------------------------------------------------------------------
# Overview sheet: just the status text + the two locale column widths.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ----------------------------------------------------------------------
# zh-cn sheet: status text, target/handback file links, wider columns.
# This locale's handback landed at 06:58:26.
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/5950f920-e0d3-4c1e-9147-b86f18c5a9cf.md", "", "", "5950f920-e0d3-4c1e-9147-b86f18c5a9cf.md") | Out-Null
$wsZh.Range("J2").Value = "5950f920-e0d3-4c1e-9147-b86f18c5a9cf.8a2d1996e1c92ae297805f6ddcd0a3b1086381cc.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-27 06:58:26"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.md", "", "", "dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.md") | Out-Null
$wsZh.Range("J3").Value = "dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.8eb2d54313dd4cd85116b56048ddb09644c904e0.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-27 06:58:26"

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(9).ColumnWidth = $maxColWidth
$wsZh.Columns.Item(10).ColumnWidth = $maxColWidth

# ----------------------------------------------------------------------
# de-de sheet: same shape, but this locale's handback landed at
# 06:58:33, so K2/K3 are written explicitly with the new timestamp.
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/5950f920-e0d3-4c1e-9147-b86f18c5a9cf.md", "", "", "5950f920-e0d3-4c1e-9147-b86f18c5a9cf.md") | Out-Null
$wsDe.Range("J2").Value = "5950f920-e0d3-4c1e-9147-b86f18c5a9cf.8a2d1996e1c92ae297805f6ddcd0a3b1086381cc.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-27 06:58:33"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.md", "", "", "dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.md") | Out-Null
$wsDe.Range("J3").Value = "dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.8eb2d54313dd4cd85116b56048ddb09644c904e0.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-27 06:58:33"

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(9).ColumnWidth = $maxColWidth
$wsDe.Columns.Item(10).ColumnWidth = $maxColWidth
